$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-10 Friday" "2024-05-11 Saturday"

Replace-Text "476÷4=" "733÷8="
Replace-Text "231÷8=" "911÷6="
Replace-Text "464÷3=" "523÷8="
Replace-Text "750÷5=" "228÷7="
Replace-Text "314÷4=" "957÷8="
Replace-Text "513÷8=" "352÷8="
Replace-Text "481÷6=" "536÷2="
Replace-Text "428÷8=" "229÷8="
Replace-Text "930÷7=" "913÷2="
Replace-Text "781÷8=" "424÷6="
Replace-Text "653÷5=" "489÷9="
Replace-Text "476÷8=" "287÷5="
Replace-Text "938÷2=" "694÷8="
Replace-Text "360÷2=" "182÷7="
Replace-Text "892÷9=" "533÷5="
Replace-Text "528÷2=" "534÷2="
Replace-Text "551÷2=" "391÷5="
Replace-Text "952÷4=" "909÷5="
Replace-Text "152÷2=" "178÷7="
Replace-Text "382÷4=" "699÷4="
Replace-Text "505÷9=" "669÷4="
Replace-Text "785÷3=" "925÷7="
Replace-Text "765÷3=" "978÷4="
Replace-Text "199÷9=" "817÷9="
Replace-Text "684÷5=" "192÷8="
